# Generate Report for Handoff
# Updates the localization status report: the file that was "c7eb5301-...md"
# has now become "Ready for handoff" (it was "In Translation"), which causes
# the rows in every sheet to be re-sorted: the "In Translation" group (still
# just 19dce665) stays first, and the "Ready for handoff" group is listed in
# alphabetical order of file name (079b75bf, 195e775a, 420445a4, c7eb5301,
# f3830def). The handoff date/time stamps for the refreshed group are bumped
# to the new handoff time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A3").Value = "079b75bf-f888-45b1-859f-74e6937ae65e.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-43-09 11:43:50"

$ov.Range("A4").Value = "195e775a-3959-4a1d-b4cb-0f7d9e55cf35.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-43-09 11:43:50"

$ov.Range("A5").Value = "420445a4-22ed-493a-aa38-69b8524c0934.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-43-09 11:43:50"

$ov.Range("A6").Value = "c7eb5301-0042-4360-b7ca-f2a861a15827.md"
$ov.Range("B6").Value = "Ready for handoff"
$ov.Range("C6").Value = "Ready for handoff"
$ov.Range("D6").Value = "2016-43-09 11:43:50"

$ov.Range("A7").Value = "f3830def-6a91-451d-9e32-682f60bab417.md"
$ov.Range("B7").Value = "Ready for handoff"
$ov.Range("C7").Value = "Ready for handoff"
$ov.Range("D7").Value = "2016-43-09 11:43:50"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A3").Value = "079b75bf-f888-45b1-859f-74e6937ae65e.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "079b75bf-f888-45b1-859f-74e6937ae65e.b581f6a80d3ddc7936153f4c68ea496ab8aaedeb.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-09 11:43:40"

$zh.Range("A4").Value = "195e775a-3959-4a1d-b4cb-0f7d9e55cf35.md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "195e775a-3959-4a1d-b4cb-0f7d9e55cf35.c25cefd95c6e1c1b5332c56d14929294f45ebc50.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-09 11:43:40"

$zh.Range("A5").Value = "420445a4-22ed-493a-aa38-69b8524c0934.md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "420445a4-22ed-493a-aa38-69b8524c0934.bad69fce416648f1a9747c041dcd3b7728a05a8b.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-09 11:43:40"

$zh.Range("A6").Value = "c7eb5301-0042-4360-b7ca-f2a861a15827.md"
$zh.Range("C6").Value = "Ready for handoff"
$zh.Range("D6").Value = "c7eb5301-0042-4360-b7ca-f2a861a15827.b956e32c1d6f6b52111057e4de75cbe8ecd2cce2.zh-cn.xlf"
$zh.Range("E6").Value = "2016-03-09 11:43:40"

$zh.Range("A7").Value = "f3830def-6a91-451d-9e32-682f60bab417.md"
$zh.Range("C7").Value = "Ready for handoff"
$zh.Range("D7").Value = "f3830def-6a91-451d-9e32-682f60bab417.21c487069c186751a8e6060e7bb06ed9b0ae5a60.zh-cn.xlf"
$zh.Range("E7").Value = "2016-03-09 11:43:40"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A3").Value = "079b75bf-f888-45b1-859f-74e6937ae65e.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "079b75bf-f888-45b1-859f-74e6937ae65e.b581f6a80d3ddc7936153f4c68ea496ab8aaedeb.de-de.xlf"
$de.Range("E3").Value = "2016-03-09 11:43:50"

$de.Range("A4").Value = "195e775a-3959-4a1d-b4cb-0f7d9e55cf35.md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "195e775a-3959-4a1d-b4cb-0f7d9e55cf35.c25cefd95c6e1c1b5332c56d14929294f45ebc50.de-de.xlf"
$de.Range("E4").Value = "2016-03-09 11:43:50"

$de.Range("A5").Value = "420445a4-22ed-493a-aa38-69b8524c0934.md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "420445a4-22ed-493a-aa38-69b8524c0934.bad69fce416648f1a9747c041dcd3b7728a05a8b.de-de.xlf"
$de.Range("E5").Value = "2016-03-09 11:43:50"

$de.Range("A6").Value = "c7eb5301-0042-4360-b7ca-f2a861a15827.md"
$de.Range("C6").Value = "Ready for handoff"
$de.Range("D6").Value = "c7eb5301-0042-4360-b7ca-f2a861a15827.b956e32c1d6f6b52111057e4de75cbe8ecd2cce2.de-de.xlf"
$de.Range("E6").Value = "2016-03-09 11:43:50"

$de.Range("A7").Value = "f3830def-6a91-451d-9e32-682f60bab417.md"
$de.Range("C7").Value = "Ready for handoff"
$de.Range("D7").Value = "f3830def-6a91-451d-9e32-682f60bab417.21c487069c186751a8e6060e7bb06ed9b0ae5a60.de-de.xlf"
$de.Range("E7").Value = "2016-03-09 11:43:50"
